$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert four new columns before the existing "ExpPoints" column (C),
# shifting it to column G. This also shifts the column styles along.
$ws.Range("C1:F1").EntireColumn.Insert()

# New header labels for the inserted columns (row 1 keeps the bold/
# centered header style that came along with the column insert).
$ws.Range("C1").Value = "WIN"
$ws.Range("D1").Value = "TOP2"
$ws.Range("E1").Value = "TOP4"
$ws.Range("F1").Value = "RELEGATION"

# Materialize empty placeholder cells in the new columns for every data
# row (2-19) so the cells exist (blank) ready for the upcoming Monte
# Carlo simulation values.
$ws.Range("C2:F19").NumberFormat = "General"
